$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) cells whose new text looks like a plain number are
# written with a leading apostrophe so Excel keeps them as text, just
# like the existing (already-text) price cells in this sheet.

$ws.Range('D2').Value = '66.831.45'
$ws.Range('E2').Value = '  +1.26%  '
$ws.Range('D3').Value = '3.283.68'
$ws.Range('E3').Value = '  -1.06%  '
$ws.Range("D4").Value = "'0.998"
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range("D5").Value = "'574.36"
$ws.Range('E5').Value = '  -1.15%  '
$ws.Range("D6").Value = "'174.84"
$ws.Range('E6').Value = '  -5.33%  '
$ws.Range("D7").Value = "'1.00"
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range("D8").Value = "'0.581"
$ws.Range('E8').Value = '  +1.43%  '
$ws.Range('D9').Value = '3.280.12'
$ws.Range('E9').Value = '  -1.04%  '
$ws.Range("D10").Value = "'0.174"
$ws.Range('E10').Value = '  -2.87%  '
$ws.Range("D11").Value = "'0.573"
$ws.Range('E11').Value = '  -0.45%  '
$ws.Range("D12").Value = "'45.41"
$ws.Range('E12').Value = '  -3.12%  '
$ws.Range('E13').Value = '  +0.35%  '
$ws.Range("D14").Value = "'688.06"
$ws.Range('E14').Value = '  +3.18%  '
$ws.Range('D15').Value = '3.810.50'
$ws.Range('E15').Value = '  -1.01%  '
$ws.Range("D16").Value = "'8.28"
$ws.Range('E16').Value = '  -1.66%  '
$ws.Range('D17').Value = '66.928.18'
$ws.Range('E18').Value = '  +1.18%  '
$ws.Range('D19').Value = '3.284.25'
$ws.Range('E19').Value = '  -1.24%  '
$ws.Range("D20").Value = "'17.28"
$ws.Range('E20').Value = '  -3.16%  '
$ws.Range("D21").Value = "'10.71"
$ws.Range('E21').Value = '  -3.09%  '
$ws.Range('E22').Value = '  -0.74%  '
$ws.Range("D23").Value = "'16.94"
$ws.Range('E23').Value = '  -4.90%  '
$ws.Range("D24").Value = "'5.16"
$ws.Range('E24').Value = '  +3.04%  '
$ws.Range("D25").Value = "'98.86"
$ws.Range('E25').Value = '  -2.19%  '
$ws.Range('E26').Value = '  -2.94%  '
$ws.Range('E27').Value = '  -2.94%  '
$ws.Range("D28").Value = "'9.23"
$ws.Range('E28').Value = '  -2.66%  '
$ws.Range("D29").Value = "'33.60"
$ws.Range('E29').Value = '  +7.33%  '
$ws.Range('E30').Value = '  -0.85%  '
$ws.Range("D31").Value = "'6.71"
$ws.Range('E31').Value = '  +0.94%  '
$ws.Range("D32").Value = "'570.30"
$ws.Range('E32').Value = '  -3.23%  '
$ws.Range('D33').Value = '3.865.16'
$ws.Range('E33').Value = '  +0.82%  '
$ws.Range("D34").Value = "'10.83"
$ws.Range('E34').Value = '  -1.19%  '
$ws.Range("D35").Value = "'0.102"
$ws.Range('E35').Value = '  -2.57%  '
$ws.Range('E36').Value = '  -0.15%  '
$ws.Range('B37').Value = 'dogwifhat'
$ws.Range('C37').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D37").Value = "'3.32"
$ws.Range('E37').Value = '  -13.40%  '
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D38").Value = "'54.97"
$ws.Range('E38').Value = '  -1.41%  '
$ws.Range('E39').Value = '  +1.40%  '
$ws.Range("D40").Value = "'3.39"
$ws.Range('E40').Value = '  -0.54%  '
$ws.Range("D41").Value = "'2.57"
$ws.Range('E41').Value = '  -3.44%  '
$ws.Range("D42").Value = "'31.71"
$ws.Range('E42').Value = '  -3.16%  '
$ws.Range('D43').Value = '0.0₃0667'
$ws.Range('E43').Value = '  -4.01%  '
$ws.Range('E44').Value = '  -2.81%  '
$ws.Range("D45").Value = "'2.97"
$ws.Range('E45').Value = '  -6.10%  '
$ws.Range('E46').Value = '  -2.03%  '
$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D47").Value = "'0.127"
$ws.Range('E47').Value = '  -0.07%  '
$ws.Range('B48').Value = 'FirstDigitalUSD'
$ws.Range('C48').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D48").Value = "'1.00"
$ws.Range('E48').Value = '  +0.07%  '
$ws.Range("D49").Value = "'2.54"
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('E50').Value = '  +5.99%  '
$ws.Range("D51").Value = "'129.91"
$ws.Range('E51').Value = '  -0.17%  '
